# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (want-to-go count) figures scraped for a handful of
# events. The same event can appear both on its category sheet (展览/演出/
# 本地生活) and on the combined "全部类型" roll-up sheet, so both copies are
# updated together.

$wb = $excel.ActiveWorkbook

$wsExhibit  = $wb.Worksheets.Item("展览")
$wsShow     = $wb.Worksheets.Item("演出")
$wsLocal    = $wb.Worksheets.Item("本地生活")
$wsAll      = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions) sheet
$wsExhibit.Range("F3").Value = 356
$wsExhibit.Range("F10").Value = 2419
$wsExhibit.Range("F17").Value = 4285
$wsExhibit.Range("F19").Value = 243

# 演出 (Shows) sheet
$wsShow.Range("F23").Value = 92

# 本地生活 (Local life) sheet
$wsLocal.Range("F4").Value = 2147

# 全部类型 (All types) roll-up sheet mirrors the rows above
$wsAll.Range("F4").Value = 2147
$wsAll.Range("F8").Value = 356
$wsAll.Range("F22").Value = 2419
$wsAll.Range("F36").Value = 243
$wsAll.Range("F48").Value = 92
